# Update countries & provincias Spain
# - Re-sorts a block of country rows (their labels shift to the
#   alphabetically-adjacent country while keeping the same row/position)
# - Refreshes the Covid-19 counters for the affected rows plus the USA (row 9)
#   and Canada (row 22)
# - Bumps the "last updated" timestamp in row 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{ Row=1;   A='Datos actualizados a 20 de Marzo de 2020 a las 00:15' }
  @{ Row=9;   A='Estados Unidos';              B=13783; C=4524; D=108; E=13468; F=64; G=57; H=207 }
  @{ Row=22;  A='Canada';                      B=873;   C=146;  D=11;  E=850;   F=1;  G=3;  H=12 }
  @{ Row=63;  A='Argentina';                   B=111;   C=14;   D=3;   E=105;   F=0;  G=1;  H=3 }
  @{ Row=64;  A='Croacia';                     B=110;   C=21;   D=5;   E=104;   F=0;  G=1;  H=1 }
  @{ Row=65;  A='Panama';                      B=109;   C=0;    D=0;   E=108;   F=7;  G=0;  H=1 }
  @{ Row=66;  A='Colombia';                    B=108;   C=15;   D=1;   E=107;   F=0;  G=0;  H=0 }
  @{ Row=67;  A='Taiwan';                      B=108;   C=8;    D=26;  E=81;    F=0;  G=0;  H=1 }
  @{ Row=68;  A='Bulgaria';                    B=107;   C=15;   D=0;   E=104;   F=0;  G=1;  H=3 }
  @{ Row=69;  A='Serbia';                      B=103;   C=14;   D=1;   E=102;   F=4;  G=0;  H=0 }
  @{ Row=103; A='Liechtenstein';               B=28;    C=0;    D=0;   E=28;    F=0;  G=0;  H=0 }
  @{ Row=104; A='Reunion';                     B=28;    C=14;   D=0;   E=28;    F=0;  G=0;  H=0 }
  @{ Row=111; A='Guayana Francesa';            B=15;    C=0;    D=0;   E=15;    F=0;  G=0;  H=0 }
  @{ Row=112; A='Bolivia';                     B=15;    C=3;    D=0;   E=15;    F=0;  G=0;  H=0 }
  @{ Row=115; A='Maldivas';                    B=13;    C=0;    D=0;   E=13;    F=0;  G=0;  H=0 }
  @{ Row=117; A='Montenegro';                  B=13;    C=5;    D=0;   E=13;    F=0;  G=0;  H=0 }
  @{ Row=118; A='Honduras';                    B=12;    C=3;    D=0;   E=12;    F=0;  G=0;  H=0 }
  @{ Row=119; A='Guam';                        B=12;    C=4;    D=0;   E=12;    F=0;  G=0;  H=0 }
  @{ Row=121; A='Ruanda';                      B=11;    C=0;    D=0;   E=11;    F=0;  G=0;  H=0 }
  @{ Row=123; A='Paraguay';                    B=11;    C=0;    D=0;   E=11;    F=1;  G=0;  H=0 }
  @{ Row=130; A='Etiopia';                     B=7;     C=1;    D=0;   E=7;     F=0;  G=0;  H=0 }
  @{ Row=131; A='Mauricio';                    B=7;     C=4;    D=0;   E=7;     F=0;  G=0;  H=0 }
  @{ Row=132; A='Kenia';                       B=7;     C=0;    D=0;   E=7;     F=0;  G=0;  H=0 }
  @{ Row=142; A='Mayotte';                     B=4;     C=1;    D=0;   E=4;     F=0;  G=0;  H=0 }
  @{ Row=143; A='San Martin (Parte Francesa)'; B=3;     C=0;    D=0;   E=3;     F=0;  G=0;  H=0 }
  @{ Row=144; A='Gabon';                       B=3;     C=0;    D=0;   E=3;     F=0;  G=0;  H=0 }
  @{ Row=145; A='Namibia';                     B=3;     C=1;    D=0;   E=3;     F=0;  G=0;  H=0 }
  @{ Row=146; A='San Bartolome';               B=3;     C=0;    D=0;   E=3;     F=0;  G=0;  H=0 }
  @{ Row=147; A='Bahamas';                     B=3;     C=2;    D=0;   E=3;     F=0;  G=0;  H=0 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.A
    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $u.B
        $ws.Cells.Item($r, 3).Value = $u.C
        $ws.Cells.Item($r, 4).Value = $u.D
        $ws.Cells.Item($r, 5).Value = $u.E
        $ws.Cells.Item($r, 6).Value = $u.F
        $ws.Cells.Item($r, 7).Value = $u.G
        $ws.Cells.Item($r, 8).Value = $u.H
    }
}
